# Update countries & provincias Spain
# - Reorder "Siria" in the country list (now ranks between Gibraltar/Guadalupe
#   in the shared-string table) with refreshed case numbers, shifting the
#   rows that used to sit below it (Guadalupe/Comoras/Guyana) down by one.
# - Swap the display order of "Groenlandia" / "Islas Malvinas".
# - Refresh total-row (row 4, "Estados Unidos") and row 20 ("Turquia") figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos ------------------------------------------------
$ws.Range("B4").Value = 2084922
$ws.Range("C4").Value = 18521
$ws.Range("D4").Value = 811922
$ws.Range("E4").Value = 1157081
$ws.Range("G4").Value = 789
$ws.Range("H4").Value = 115919

# --- Row 20: Turquia -------------------------------------------------------
$ws.Range("B20").Value = 97476
$ws.Range("C20").Value = 351
$ws.Range("D20").Value = 57608
$ws.Range("E20").Value = 31872

# --- Rows 166-169: Siria moves up (new data), Guadalupe/Comoras/Guyana shift
#     down one row keeping their previous figures ---------------------------
$ws.Range("A166").Value = "Siria"
$ws.Range("B166").Value = 164
$ws.Range("C166").Value = 12
$ws.Range("D166").Value = 68
$ws.Range("E166").Value = 90
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 6

$ws.Range("A167").Value = "Guadalupe"
$ws.Range("B167").Value = 164
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 144
$ws.Range("E167").Value = 6
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 14

$ws.Range("A168").Value = "Comoras"
$ws.Range("B168").Value = 162
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 97
$ws.Range("E168").Value = 63
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 2

$ws.Range("A169").Value = "Guyana"
$ws.Range("B169").Value = 156
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 92
$ws.Range("E169").Value = 52
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 12

# --- Rows 206-207: swap Groenlandia / Islas Malvinas ------------------------
$ws.Range("A206").Value = "Islas Malvinas"
$ws.Range("A207").Value = "Groenlandia"
